$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2-17).
# Every one of these needs to move forward by one day: 45179 -> 45180.
for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
